$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F7").Value = 65
$ws.Range("F8").Value = 1561
$ws.Range("F10").Value = 3392
$ws.Range("F11").Value = 466
$ws.Range("F17").Value = 1992
$ws.Range("F19").Value = 617592
$ws.Range("F22").Value = 1190
$ws.Range("F26").Value = 1917
$ws.Range("F29").Value = 610
$ws.Range("F35").Value = 1920
$ws.Range("F36").Value = 1146
$ws.Range("F38").Value = 128
$ws.Range("F45").Value = 2947
$ws = $wb.Worksheets.Item(2)
$ws.Range("F13").Value = 136005
$ws.Range("F14").Value = 136005
$ws.Range("F29").Value = 333
$ws.Range("F41").Value = 171
$ws = $wb.Worksheets.Item(3)
$ws.Range("F9").Value = 754
$ws.Range("F10").Value = 1051
$ws.Range("F13").Value = 1492
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 754
$ws.Range("F7").Value = 1492
$ws.Range("F11").Value = 65
$ws.Range("F12").Value = 1561
$ws.Range("F14").Value = 466
$ws.Range("F19").Value = 1992
$ws.Range("F21").Value = 617594
$ws.Range("F25").Value = 136006
$ws.Range("F27").Value = 1190
$ws.Range("F31").Value = 1917
$ws.Range("F34").Value = 610
$ws.Range("F42").Value = 1146
$ws.Range("F43").Value = 128
$ws.Range("F50").Value = 2947
